# Course-Tracking-Sheet.xlsx edit:
#   - Row 30 (E30): append a blue "Arrow Function" run to the let/var/const topic text.
#   - Row 31 (D31/E31): "Constructors & Prototypes" topic gets a "this keyword, ..." entry.
#   - Row 32 (A32/D32/E32): new "Prototypes" session with "prototype, class" topic,
#     highlighted to match the surrounding rows.
#   - E18: font color tweaked to the blue Calibri used elsewhere in the sheet.
#   - Selection / active cell moves to E32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E18: recolor existing text (adds the new blue font + cellXf pair) ---
$ws.Range("E18").Font.Color = 15773696

# --- D31: "Constructors & Prototypes" moves down into D31 ---
$ws.Range("D31").Value = "Constructors & Prototypes"

# --- D32: new "Prototypes" entry, yellow-highlighted like its neighbours ---
$ws.Range("D32").Value = "Prototypes"
$ws.Range("D32").Interior.Color = 65535

# --- E30: "let, var, const, self invoking function, closure" + blue "Arrow Function" ---
$e30 = $ws.Range("E30")
$prefix = "let, var, const, self invoking function, closure, "
$suffix = "Arrow Function"
$e30.Value = $prefix + $suffix
$startPos = $prefix.Length + 1
$chars = $e30.Characters($startPos, $suffix.Length)
$chars.Font.Color = 15773696

# --- E31: this keyword, constuctor function, call, apply, bind ---
$ws.Range("E31").Value = "this keyword, constuctor function, call, apply, bind"

# --- A32: blue highlight matching the rest of the block ---
$ws.Range("A32").Interior.Color = 15773696

# --- E32: topics for the new Prototypes row ---
$ws.Range("E32").Value = "prototype, class"

# --- Selection moves to E32 (matches the author's final cursor position) ---
$ws.Range("E32").Select() | Out-Null

Write-Output "Edit applied"
